$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose 7-day window no longer has a full trailing sample (need to become blank)
$blankRows = @(5, 6, 7)
foreach ($r in $blankRows) {
    $ws.Cells.Item($r, 3).Value = "'"
    $ws.Cells.Item($r, 3).Style = "Normal"
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Cells.Item($r, 4).Style = "Normal"
}

# Row -> [C value, D value] for the recentred (trailing) 7-day window
$data = @{
    8 = @(1, 7.852375343541421)
    12 = @(2, 15.70475068708284)
    13 = @(2, 15.70475068708284)
    14 = @(2, 15.70475068708284)
    15 = @(4, 31.40950137416569)
    16 = @(3, 23.55712603062426)
    19 = @(4, 31.40950137416569)
    20 = @(3, 23.55712603062426)
    21 = @(3, 23.55712603062426)
    22 = @(1, 7.852375343541421)
    25 = @(2, 15.70475068708284)
    29 = @(1, 7.852375343541421)
    30 = @(1, 7.852375343541421)
    31 = @(1, 7.852375343541421)
    33 = @(0, 0)
    34 = @(0, 0)
    35 = @(0, 0)
    36 = @(1, 7.852375343541421)
    37 = @(1, 7.852375343541421)
    38 = @(1, 7.852375343541421)
    40 = @(2, 15.70475068708284)
    41 = @(2, 15.70475068708284)
    43 = @(1, 7.852375343541421)
    44 = @(1, 7.852375343541421)
    45 = @(2, 15.70475068708284)
    47 = @(3, 23.55712603062426)
    48 = @(3, 23.55712603062426)
    49 = @(3, 23.55712603062426)
    51 = @(7, 54.96662740478995)
    53 = @(4, 31.40950137416569)
    54 = @(4, 31.40950137416569)
    55 = @(6, 47.11425206124853)
    56 = @(8, 62.81900274833137)
    57 = @(8, 62.81900274833137)
    58 = @(17, 133.4903808402042)
    59 = @(21, 164.8998822143699)
    60 = @(22, 172.7522575579113)
    61 = @(28, 219.8665096191598)
    62 = @(34, 266.9807616804083)
    63 = @(32, 251.2760109933255)
    65 = @(31, 243.423635649784)
    66 = @(35, 274.8331370239497)
    67 = @(37, 290.5378877110326)
    68 = @(35, 274.8331370239497)
    69 = @(31, 243.423635649784)
    70 = @(40, 314.0950137416569)
    71 = @(45, 353.356890459364)
    72 = @(52, 408.3235178641539)
    73 = @(54, 424.0282685512367)
    74 = @(60, 471.1425206124852)
    75 = @(64, 502.552021986651)
    76 = @(68, 533.9615233608166)
    77 = @(66, 518.2567726737337)
    78 = @(72, 565.3710247349823)
    79 = @(63, 494.6996466431095)
    80 = @(59, 463.2901452689439)
    81 = @(56, 439.7330192383196)
    82 = @(51, 400.4711425206125)
    83 = @(55, 431.8806438947781)
    84 = @(55, 431.8806438947781)
    85 = @(43, 337.6521397722811)
    86 = @(47, 369.0616411464468)
    87 = @(46, 361.2092658029054)
    88 = @(40, 314.0950137416569)
    89 = @(44, 345.5045151158225)
    90 = @(36, 282.6855123674911)
    91 = @(36, 282.6855123674911)
    92 = @(34, 266.9807616804083)
    93 = @(28, 219.8665096191598)
    94 = @(26, 204.161758932077)
    95 = @(33, 259.1283863368669)
    96 = @(34, 266.9807616804083)
    97 = @(39, 306.2426383981154)
    98 = @(34, 266.9807616804083)
    99 = @(38, 298.390263054574)
    100 = @(43, 337.6521397722811)
    101 = @(41, 321.9473890851983)
    102 = @(36, 282.6855123674911)
    103 = @(29, 227.7188849627012)
    104 = @(28, 219.8665096191598)
    105 = @(38, 298.390263054574)
    106 = @(32, 251.2760109933255)
    107 = @(27, 212.0141342756183)
    108 = @(27, 212.0141342756183)
    109 = @(27, 212.0141342756183)
    110 = @(35, 274.8331370239497)
    111 = @(31, 243.423635649784)
    112 = @(22, 172.7522575579113)
    113 = @(26, 204.161758932077)
    114 = @(26, 204.161758932077)
    115 = @(26, 204.161758932077)
    117 = @(28, 219.8665096191598)
    118 = @(28, 219.8665096191598)
    119 = @(29, 227.7188849627012)
    120 = @(24, 188.4570082449941)
    121 = @(24, 188.4570082449941)
    123 = @(29, 227.7188849627012)
    124 = @(26, 204.161758932077)
    125 = @(25, 196.3093835885355)
    126 = @(27, 212.0141342756183)
    127 = @(28, 219.8665096191598)
    129 = @(34, 266.9807616804083)
    130 = @(32, 251.2760109933255)
    131 = @(28, 219.8665096191598)
    132 = @(31, 243.423635649784)
    133 = @(28, 219.8665096191598)
    134 = @(34, 266.9807616804083)
    135 = @(33, 259.1283863368669)
    136 = @(30, 235.5712603062426)
    137 = @(29, 227.7188849627012)
    138 = @(30, 235.5712603062426)
    139 = @(27, 212.0141342756183)
    140 = @(24, 188.4570082449941)
    141 = @(20, 157.0475068708284)
    142 = @(19, 149.195131527287)
    143 = @(19, 149.195131527287)
    144 = @(24, 188.4570082449941)
    145 = @(28, 219.8665096191598)
    146 = @(38, 298.390263054574)
    147 = @(41, 321.9473890851983)
    148 = @(40, 314.0950137416569)
    149 = @(40, 314.0950137416569)
    150 = @(40, 314.0950137416569)
    151 = @(36, 282.6855123674911)
    152 = @(34, 266.9807616804083)
    153 = @(22, 172.7522575579113)
    154 = @(25, 196.3093835885355)
    156 = @(26, 204.161758932077)
    157 = @(27, 212.0141342756183)
    158 = @(24, 188.4570082449941)
    159 = @(24, 188.4570082449941)
    161 = @(21, 164.8998822143699)
    162 = @(26, 204.161758932077)
    163 = @(25, 196.3093835885355)
    164 = @(24, 188.4570082449941)
    165 = @(25, 196.3093835885355)
    166 = @(26, 204.161758932077)
    167 = @(29, 227.7188849627012)
    168 = @(35, 274.8331370239497)
    169 = @(29, 227.7188849627012)
    170 = @(28, 219.8665096191598)
    171 = @(28, 219.8665096191598)
    172 = @(33, 259.1283863368669)
    173 = @(33, 259.1283863368669)
    174 = @(36, 282.6855123674911)
    175 = @(44, 345.5045151158225)
    176 = @(53, 416.1758932076954)
    177 = @(57, 447.585394581861)
    178 = @(63, 494.6996466431095)
    179 = @(64, 502.552021986651)
    180 = @(68, 533.9615233608166)
    181 = @(71, 557.5186493914408)
    182 = @(75, 588.9281507656066)
    183 = @(77, 604.6329014526895)
    184 = @(77, 604.6329014526895)
}
foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]
    $ws.Cells.Item($r, 4).Value = $vals[1]
}
